$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5: fill in previously-missing X5 / Y5 ---
$ws.Range("X5").Value = 0.18999999999999773
$ws.Range("Y5").Value = "Up"

# --- New data rows 6-11 ---
$rows = @(
    @{ Row=6;  A=42650.338379629633; B=11; C="Buy";     D=50; E=8137;  F=892;  G=68; H=30; I=86; J=13; K=5412;  L=149; M=67;  N=59; O=9;  P="Bag"; Q=38.48959524716075; R=0; S=0.1046; T=0.0345; U=4.82; V=2.2799999999999998; W=0; X=0.18999999999999773;  Y="Up"   }
    @{ Row=7;  A=42650.339618055557; B=-5; C="Neutral"; D=0;  E=4;     F=2;    G=0;  H=0;  I=0;  J=0;  K=56;    L=0;   M=0;   N=0;  O=0;  P="Bag"; Q=38.48959524716075; R=0; S=0.1046; T=0.0345; U=4.82; V=2.2799999999999998; W=0; X=0.18999999999999773;  Y="Up"   }
    @{ Row=8;  A=42650.348773148151; B=11; C="Buy";     D=32; E=33453; F=3627; G=61; H=37; I=86; J=13; K=17329; L=514; M=314; N=84; O=13; P="Bag"; Q=38.48959524716075; R=0; S=0.1046; T=0.0345; U=4.82; V=2.2799999999999998; W=0; X=0.18999999999999773;  Y="Up"   }
    @{ Row=9;  A=42650.359050925923; B=1;  C="Neutral"; D=2;  E=2806;  F=315;  G=57; H=42; I=50; J=50; K=2767;  L=41;  M=30;  N=1;  O=1;  P="Bag"; Q=38.48959524716075; R=0; S=0.1046; T=0.0345; U=4.82; V=2.2799999999999998; W=0; X=0.18999999999999773;  Y="Up"   }
    @{ Row=10; A=42650.36146990741;  B=11; C="Buy";     D=50; E=8149;  F=893;  G=68; H=30; I=86; J=13; K=5592;  L=149; M=67;  N=59; O=9;  P="Bag"; Q=38.48959524716075; R=0; S=0.1046; T=0.0345; U=4.82; V=2.2799999999999998; W=0; X=0.18999999999999773;  Y="Up"   }
    @{ Row=11; A=42650.36310185185;  B=1;  C="Neutral"; D=0;  E=4;     F=2;    G=0;  H=0;  I=0;  J=0;  K=11;    L=0;   M=0;   N=0;  O=0;  P="Bag"; Q=37.799019424898844; R=0; S=0.1046; T=0.0343; U=4.82; V=2.2799999999999998; W=0 }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Cells.Item($n, 1).Value = $r.A
    $ws.Cells.Item($n, 2).Value = $r.B
    $ws.Cells.Item($n, 3).Value = $r.C
    $ws.Cells.Item($n, 4).Value = $r.D
    $ws.Cells.Item($n, 5).Value = $r.E
    $ws.Cells.Item($n, 6).Value = $r.F
    $ws.Cells.Item($n, 7).Value = $r.G
    $ws.Cells.Item($n, 8).Value = $r.H
    $ws.Cells.Item($n, 9).Value = $r.I
    $ws.Cells.Item($n, 10).Value = $r.J
    $ws.Cells.Item($n, 11).Value = $r.K
    $ws.Cells.Item($n, 12).Value = $r.L
    $ws.Cells.Item($n, 13).Value = $r.M
    $ws.Cells.Item($n, 14).Value = $r.N
    $ws.Cells.Item($n, 15).Value = $r.O
    $ws.Cells.Item($n, 16).Value = $r.P
    $ws.Cells.Item($n, 17).Value = $r.Q
    $ws.Cells.Item($n, 18).Value = $r.R

    $ws.Cells.Item($n, 19).Value = $r.S
    $ws.Cells.Item($n, 19).NumberFormat = "0.00%"
    $ws.Cells.Item($n, 20).Value = $r.T
    $ws.Cells.Item($n, 20).NumberFormat = "0.00%"

    $ws.Cells.Item($n, 21).Value = $r.U
    $ws.Cells.Item($n, 22).Value = $r.V
    $ws.Cells.Item($n, 23).Value = $r.W

    if ($r.ContainsKey("X")) {
        $ws.Cells.Item($n, 24).Value = $r.X
    }
    if ($r.ContainsKey("Y")) {
        $ws.Cells.Item($n, 25).Value = $r.Y
    }
}

# --- Selection reflects the user's last click while editing ---
$ws.Range("B7").Select() | Out-Null
